$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.943.44'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '2.919.00'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.49'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.93'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.58'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").Value = '3.402.16'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").Value = '60.891.16'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.70'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").Value = '2.920.67'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '432.50'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.27'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.92'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.20'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("E26").Value = '  -0.86%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +4.80%  '
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("E31").Value = '  +3.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.66'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '0.0₃0862'
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("E35").Value = '  -0.70%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("E39").Value = '  -5.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.54'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '41.43'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  -4.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '376.25'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").Value = '2.705.52'
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.72'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("E48").Value = '  -3.86%  '
$ws.Range("E49").Value = '  -0.48%  '
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("E51").Value = '  -0.91%  '
